{"js": "// The document originally has three short paragraphs about \"Surface Area\" /\n// \"Area\" (paragraphs[0..2]), followed by a paragraph about the \"SI Unit of\n// surface area\" (paragraphs[3]), two spacer paragraphs with paragraph\n// formatting (paragraphs[4..5]), and two bare empty paragraphs\n// (paragraphs[6..7]) before the final section properties.\n//\n// The edit:\n//   1. Rewrites the first three paragraphs as a single paragraph (keeping\n//      the first paragraph's paragraph formatting) made of three runs that\n//      talk about \"Volume\" instead of \"Surface area\" / \"Area\", with\n//      grammar proofing markers (<w:proofErr/>) bracketing the middle run\n//      (\"It's\").\n//   2. Removes the \"SI Unit of surface area...\" paragraph, the two spacer\n//      paragraphs, and one of the two trailing bare paragraphs, leaving a\n//      single trailing empty paragraph before the section break.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// ---- Step 1: merge paragraphs[0..2] into a single paragraph with new text ----\n\nconst mergeStart = paragraphs.items[0].getRange(\"Start\");\n// Use the start of paragraphs[3] (\"SI Unit of...\") as the (exclusive) end\n// boundary so the expanded range swallows paragraphs[0..2] *including*\n// their paragraph marks.\nconst mergeEndBoundary = paragraphs.items[3].getRange(\"Start\");\nconst mergeRange = mergeStart.expandTo(mergeEndBoundary);\n\nconst runProps = '<w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:eastAsia=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:color w:val=\"222222\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr>';\nconst paraProps = '<w:pPr><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/><w:spacing w:before=\"120\" w:after=\"120\" w:line=\"240\" w:lineRule=\"auto\"/>' + runProps + '</w:pPr>';\n\nconst run1Text = 'Volume is quantity of three dimensional space enclosed by closed surface, the space that it contains. It is derived form of length. ';\nconst run2Text = 'It\\u2019s';\nconst run3Text = ' unit is in cubic meter. Some have regular shape and it has formula. Volume is calculated by integration formula.';\n\nconst newParaXml =\n  '<w:p>' + paraProps +\n  '<w:r>' + runProps + '<w:t xml:space=\"preserve\">' + run1Text + '</w:t></w:r>' +\n  '<w:proofErr w:type=\"gramStart\"/>' +\n  '<w:r>' + runProps + '<w:t>' + run2Text + '</w:t></w:r>' +\n  '<w:proofErr w:type=\"gramEnd\"/>' +\n  '<w:r>' + runProps + '<w:t xml:space=\"preserve\">' + run3Text + '</w:t></w:r>' +\n  '</w:p>';\n\nconst packageXml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + newParaXml + '</w:body></w:document></pkg:xmlData>' +\n  '</pkg:part></pkg:package>';\n\nmergeRange.insertOoxml(packageXml, Word.InsertLocation.replace);\nawait context.sync();\n\n// ---- Step 2: delete the \"SI Unit of...\" paragraph, the two spacer\n//      paragraphs, and one trailing bare paragraph ----\n\nconst paragraphsAfterMerge = context.document.body.paragraphs;\nparagraphsAfterMerge.load(\"items\");\nawait context.sync();\n\n// After the merge, index 0 is the new \"Volume...\" paragraph, index 1 is\n// \"SI Unit of surface area...\", indices 2-3 are the spacer paragraphs, and\n// indices 4-5 are the two bare trailing paragraphs. Delete through the\n// start of index 5 so one bare paragraph remains.\nconst deleteStart = paragraphsAfterMerge.items[1].getRange(\"Start\");\nconst deleteEndBoundary = paragraphsAfterMerge.items[5].getRange(\"Start\");\nconst deleteRange = deleteStart.expandTo(deleteEndBoundary);\ndeleteRange.delete();\nawait context.sync();\n", "ps1": "# The document originally has three short paragraphs about \"Surface Area\" /\n# \"Area\" (paragraphs 1-3), followed by a paragraph about the \"SI Unit of\n# surface area\" (paragraph 4), two spacer paragraphs with paragraph\n# formatting (paragraphs 5-6), and two bare empty paragraphs (paragraphs\n# 7-8) before the final section properties.\n#\n# The edit:\n#   1. Rewrites the first three paragraphs as a single paragraph (keeping\n#      paragraph 1's paragraph formatting) made of three runs that talk\n#      about \"Volume\" instead of \"Surface area\" / \"Area\", with grammar\n#      proofing markers (<w:proofErr/>) bracketing the middle run (\"It's\").\n#   2. Removes the \"SI Unit of surface area...\" paragraph, the two spacer\n#      paragraphs, and one of the two trailing bare paragraphs, leaving a\n#      single trailing empty paragraph before the section break.\n\n$d = $word.ActiveDocument\n\n# ---- Step 1: merge paragraphs 1-3 into a single paragraph with new text ----\n\n$firstPara = $d.Paragraphs(1)\n$thirdPara = $d.Paragraphs(3)\n$mergeRange = $d.Range($firstPara.Range.Start, $thirdPara.Range.End)\n\n$runProps = '<w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:eastAsia=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:color w:val=\"222222\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr>'\n$paraProps = '<w:pPr><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/><w:spacing w:before=\"120\" w:after=\"120\" w:line=\"240\" w:lineRule=\"auto\"/>' + $runProps + '</w:pPr>'\n\n$apostrophe = [char]0x2019\n$run1Text = 'Volume is quantity of three dimensional space enclosed by closed surface, the space that it contains. It is derived form of length. '\n$run2Text = 'It' + $apostrophe + 's'\n$run3Text = ' unit is in cubic meter. Some have regular shape and it has formula. Volume is calculated by integration formula.'\n\n$newParaXml = '<w:p>' + $paraProps + `\n    '<w:r>' + $runProps + '<w:t xml:space=\"preserve\">' + $run1Text + '</w:t></w:r>' + `\n    '<w:proofErr w:type=\"gramStart\"/>' + `\n    '<w:r>' + $runProps + '<w:t>' + $run2Text + '</w:t></w:r>' + `\n    '<w:proofErr w:type=\"gramEnd\"/>' + `\n    '<w:r>' + $runProps + '<w:t xml:space=\"preserve\">' + $run3Text + '</w:t></w:r>' + `\n    '</w:p>'\n\n$packageXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n    '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData>' + `\n    '</pkg:part></pkg:package>'\n\n$mergeRange.InsertXML($packageXml)\n\n# ---- Step 2: delete the \"SI Unit of...\" paragraph, the two spacer\n#      paragraphs, and one trailing bare paragraph ----\n\n$siUnitPara = $d.Paragraphs(2)\n$lastBarePara = $d.Paragraphs(5)\n$deleteRange = $d.Range($siUnitPara.Range.Start, $lastBarePara.Range.End)\n$deleteRange.Delete()\n"}
